$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string text cells (E column) to render markdown line breaks
# using real newlines instead of literal "<br>" text.
$newlineChar = [char]10

$lowText   = "Supports the GCS if coverage is **Low** " + $newlineChar + "Other members: Global South + EU " + $newlineChar + "(25-33% of world emissions)"
$midText   = "Supports the GCS if coverage is **Mid** " + $newlineChar + "Global South + China " + $newlineChar + "(56% of world emissions)"
$highText  = "Supports the GCS if coverage is **High** " + $newlineChar + "Global South + China + EU + various HICs " + $newlineChar + "(UK, Japan, Korea, Canada...; 64-72% of emissions)"
$colorText = "Supports the GCS if coverage is **High**, **color** variant " + $newlineChar + "Global South + China + EU + various HICs " + $newlineChar + "+ Distributive effects shown using colors on world map"

# Find the E-column cells referencing each text (rows 2..73) and replace them.
$usedRange = $ws.UsedRange
$maxRow = $usedRange.Rows.Count

for ($r = 2; $r -le $maxRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value()
    if ($val -eq $null) { continue }
    if ($val -like "Supports the GCS if coverage is **High**, **color** variant<br>*") {
        $cell.Value = $colorText
    } elseif ($val -like "Supports the GCS if coverage is **Low**<br>*") {
        $cell.Value = $lowText
    } elseif ($val -like "Supports the GCS if coverage is **Mid**<br>*") {
        $cell.Value = $midText
    } elseif ($val -like "Supports the GCS if coverage is **High**<br>*") {
        $cell.Value = $highText
    }
}

# --- Update the slightly-adjusted mean / CI_low / CI_high values for the
# specific rows identified in the diff.
$rowUpdates = @{
    2  = @(67.6970977479236, 66.3023963959884, 69.0917990998588)
    12 = @(73.6179977793817, 69.6993485561823, 77.536647002581)
    14 = @(57.1528754247366, 56.2673732215499, 58.0383776279233)
    24 = @(49.0525173251567, 45.9466877611121, 52.1583468892013)
    38 = @(67.0197929692514, 65.3075604359723, 68.7320255025304)
    48 = @(64.3326161149391, 58.5319321064818, 70.1333001233965)
    50 = @(69.9117027281215, 68.2941286055139, 71.529276850729)
    60 = @(59.5603805397135, 53.586463176426, 65.5342979030011)
    62 = @(63.5326678113447, 61.8124249135943, 65.2529107090951)
    72 = @(53.7266999164565, 47.5942502637921, 59.8591495691209)
}

foreach ($row in $rowUpdates.Keys) {
    $vals = $rowUpdates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
